# Auto-generated: apply scheduled market-price refresh to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1086.2222
$ws.Range("I18").Value = 988.3077
$ws.Range("K18").Value = 988.3077
$ws.Range("M18").Value = -704.3077
$ws.Range("H129").Value = 880.84784
$ws.Range("I129").Value = 334.66666
$ws.Range("J129").Value = 918.9535
$ws.Range("K129").Value = 1003.99998
$ws.Range("L129").Value = 2756.8605
$ws.Range("M129").Value = 3996.00002
$ws.Range("N129").Value = -12756.8605
$ws.Range("H131").Value = 995
$ws.Range("I131").Value = 995
$ws.Range("K131").Value = 2985
$ws.Range("M131").Value = 2055
$ws.Range("H137").Value = 1442.12
$ws.Range("I137").Value = 970.46155
$ws.Range("K137").Value = 2911.38465
$ws.Range("M137").Value = -361.38465
$ws.Range("H138").Value = 1475.3334
$ws.Range("I138").Value = 1131.0605
$ws.Range("J138").Value = 1745.8334
$ws.Range("K138").Value = 3393.1815
$ws.Range("L138").Value = 5237.5002
$ws.Range("M138").Value = 1746.8185
$ws.Range("N138").Value = -15517.5002
$ws.Range("H141").Value = 627.1429000000001
$ws.Range("I141").Value = 629.2308
$ws.Range("J141").Value = 600
$ws.Range("K141").Value = 1887.6924
$ws.Range("L141").Value = 1800
$ws.Range("M141").Value = 3292.3076
$ws.Range("N141").Value = -12160

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7313.933
$ws.Range("I2").Value = 633
$ws.Range("J2").Value = 17335.334
$ws.Range("K2").Value = 633
$ws.Range("L2").Value = 17335.334
$ws.Range("M2").Value = -520
$ws.Range("N2").Value = -17561.334
$ws.Range("H31").Value = 4718.875
$ws.Range("I31").Value = 4718.875
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4718.875
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4424.875
$ws.Range("N31").Value = $null
$ws.Range("H61").Value = 2270.6667
$ws.Range("I61").Value = 1906
$ws.Range("K61").Value = 1906
$ws.Range("M61").Value = -1694
$ws.Range("H116").Value = 7313.933
$ws.Range("I116").Value = 633
$ws.Range("J116").Value = 17335.334
$ws.Range("K116").Value = 633
$ws.Range("L116").Value = 17335.334
$ws.Range("M116").Value = 1661
$ws.Range("N116").Value = -21923.334
$ws.Range("H136").Value = 2270.6667
$ws.Range("I136").Value = 1906
$ws.Range("K136").Value = 5718
$ws.Range("M136").Value = -3168

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7313.933
$ws.Range("I3").Value = 633
$ws.Range("J3").Value = 17335.334
$ws.Range("K3").Value = 633
$ws.Range("L3").Value = 17335.334
$ws.Range("M3").Value = -519
$ws.Range("N3").Value = -17563.334
$ws.Range("H94").Value = 25001048
$ws.Range("I94").Value = 35715156
$ws.Range("J94").Value = 1459.6666
$ws.Range("K94").Value = 35715156
$ws.Range("L94").Value = 1459.6666
$ws.Range("M94").Value = -35714705
$ws.Range("N94").Value = -2361.6666
$ws.Range("H134").Value = 20924.834
$ws.Range("I134").Value = 1874.5
$ws.Range("J134").Value = 30450
$ws.Range("K134").Value = 5623.5
$ws.Range("L134").Value = 91350
$ws.Range("M134").Value = -3088.5
$ws.Range("N134").Value = -96420

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2470
$ws.Range("H125").Value = 18000
$ws.Range("J125").Value = 18000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -22920
$ws.Range("H132").Value = 11311.167
$ws.Range("I132").Value = 18622.666
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 55867.99800000001
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -53337.99800000001
$ws.Range("N132").Value = -17059.0001
$ws.Range("H134").Value = 6002.75
$ws.Range("I134").Value = 7337
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 22011
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -19476
$ws.Range("N134").Value = -11070

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 487923.84
$ws.Range("I4").Value = 179963.8
$ws.Range("J4").Value = 568965.9399999999
$ws.Range("K4").Value = 539891.3999999999
$ws.Range("L4").Value = 1706897.82
$ws.Range("M4").Value = -539779.3999999999
$ws.Range("N4").Value = -1707121.82
$ws.Range("H5").Value = 1274.9032
$ws.Range("I5").Value = 1377.8077
$ws.Range("J5").Value = 739.8
$ws.Range("K5").Value = 4133.4231
$ws.Range("L5").Value = 2219.4
$ws.Range("M5").Value = -4021.4231
$ws.Range("N5").Value = -2443.4
$ws.Range("H131").Value = 1701.97
$ws.Range("I131").Value = 320
$ws.Range("J131").Value = 1744.7113
$ws.Range("K131").Value = 960
$ws.Range("L131").Value = 5234.1339
$ws.Range("M131").Value = 4080
$ws.Range("N131").Value = -15314.1339
$ws.Range("H135").Value = 1274.9032
$ws.Range("I135").Value = 1377.8077
$ws.Range("J135").Value = 739.8
$ws.Range("K135").Value = 12400.2693
$ws.Range("L135").Value = 6658.2
$ws.Range("M135").Value = -9865.2693
$ws.Range("N135").Value = -11728.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 911
$ws.Range("I97").Value = 901.1111
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 901.1111
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -405.1111
$ws.Range("N97").Value = -1992
$ws.Range("H132").Value = 3043.75
$ws.Range("I132").Value = 3131.1667
$ws.Range("J132").Value = 3006.2856
$ws.Range("K132").Value = 9393.500100000001
$ws.Range("L132").Value = 9018.856800000001
$ws.Range("M132").Value = -6863.500100000001
$ws.Range("N132").Value = -14078.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 375573.4
$ws.Range("J2").Value = 329489.47
$ws.Range("L2").Value = 329489.47
$ws.Range("N2").Value = -329713.47
$ws.Range("H22").Value = 1663.6666
$ws.Range("I22").Value = 991
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 991
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -696
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1663.6666
$ws.Range("I27").Value = 991
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 991
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -884
$ws.Range("N27").Value = -2214
$ws.Range("J93").Value = 1000
$ws.Range("L93").Value = 1000
$ws.Range("N93").Value = -3496
$ws.Range("H125").Value = 39589.668
$ws.Range("J125").Value = 39589.668
$ws.Range("L125").Value = 39589.668
$ws.Range("N125").Value = -49429.668
$ws.Range("H132").Value = 24195.934
$ws.Range("J132").Value = 45891.39
$ws.Range("L132").Value = 137674.17
$ws.Range("N132").Value = -142734.17

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 333334000
$ws.Range("I2").Value = 500000000
$ws.Range("K2").Value = 500000000
$ws.Range("M2").Value = -499999888
$ws.Range("H31").Value = 29500
$ws.Range("J31").Value = 29500
$ws.Range("L31").Value = 29500
$ws.Range("N31").Value = -30196
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -87
$ws.Range("N39").Value = $null
$ws.Range("H132").Value = 6605.5
$ws.Range("I132").Value = 10773
$ws.Range("J132").Value = 3364.111
$ws.Range("K132").Value = 32319
$ws.Range("L132").Value = 10092.333
$ws.Range("M132").Value = -29789
$ws.Range("N132").Value = -15152.333
